$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.186737775649103

# Row 3
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("E3").Value = "Jewel"
$ws.Range("G3").Value = 6.491841461046875
$ws.Range("H3").Value = "Black or African American"

# Row 4
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("E4").Value = "Colleen"
$ws.Range("G4").Value = 6.098775820819761
$ws.Range("H4").Value = "White"

# Row 5
$ws.Range("G5").Value = 5.03843188797754

# Row 6
$ws.Range("G6").Value = 5.020761003118488

# Row 7
$ws.Range("G7").Value = 4.239120605821088

# Row 8
$ws.Range("G8").Value = 1.248319253184411

# Row 9
$ws.Range("G9").Value = 1.073514296423548

# Row 10
$ws.Range("C10").Value = 32
$ws.Range("D10").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("E10").Value = "Kellie"
$ws.Range("G10").Value = 0.4217658868825644
$ws.Range("H10").Value = "White"

# Row 11
$ws.Range("C11").Value = 21
$ws.Range("D11").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("E11").Value = "Bri"
$ws.Range("G11").Value = 0.4004490700212808

# Row 12
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("E12").Value = "Shadaisia"
$ws.Range("G12").Value = 0.3953773248513697
$ws.Range("H12").Value = "Black or African American"

# Row 13
$ws.Range("C13").Value = 33
$ws.Range("D13").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("E13").Value = "Shaniek"
$ws.Range("G13").Value = 0.1805031549777598

# Row 14
$ws.Range("G14").Value = 13.04325402792447

# Row 15
$ws.Range("G15").Value = 8.342377812971202

# Row 16
$ws.Range("G16").Value = 7.489472321657063

# Row 17
$ws.Range("G17").Value = 7.219432926815826

# Row 18
$ws.Range("C18").Value = 22
$ws.Range("D18").Value = "60db4fde6193c50664c9c478"
$ws.Range("E18").Value = "Edosagbe"
$ws.Range("G18").Value = 5.404956080902719
$ws.Range("H18").Value = "Black or African American"

# Row 19
$ws.Range("C19").Value = 32
$ws.Range("D19").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("E19").Value = "Jamarii"
$ws.Range("G19").Value = 5.203546488046102

# Row 20
$ws.Range("C20").Value = 26
$ws.Range("D20").Value = "5dd671942b033b5ec8bc97b4"
$ws.Range("E20").Value = "Juan"
$ws.Range("G20").Value = 5.194694186643499
$ws.Range("H20").Value = "Hispanic"

# Row 21
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = "60b322994d0b901954690036"
$ws.Range("E21").Value = "Brennan"
$ws.Range("G21").Value = 4.334666484926464

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = "5e2522d6b734b47915f88275"
$ws.Range("E22").Value = "Corey"
$ws.Range("G22").Value = 4.178693876440433

# Row 23
$ws.Range("G23").Value = 3.419194189605884

# Row 24
$ws.Range("G24").Value = 2.385885516067507

# Row 25
$ws.Range("G25").Value = 2.223286854337817
